$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: clear "Đơn vị:" from S6 and remove "Triệu USD" from T6 entirely
$ws.Range("S6").ClearContents()
$ws.Range("T6").Clear()

# Row 7 header labels: "Doanh số ..." -> "Hồ sơ ..."
$ws.Range("Q7").Value = "Hồ sơ chi quầy"
$ws.Range("R7").Value = "Hồ sơ chi nhà"
$ws.Range("S7").Value = "Hồ sơ chuyển khoản"
$ws.Range("T7").Value = "Tổng"

# Row 35 header labels: lowercase -> capitalized / "Tổng" unchanged
$ws.Range("Q35").Value = "Chi quầy"
$ws.Range("R35").Value = "Chi nhà"
$ws.Range("S35").Value = "Chi CK"
$ws.Range("T35").Value = "Tổng"

# View/format cosmetic changes: move top-left scroll position and selection
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("S35").Select()
